$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-off timestamp value in A4
$ws.Range("A4").Value = 45877.12518302083

# Append new row 5 with the new reading
$ws.Range("A5").Value = 45877.33355579317
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 32
$ws.Range("D5").Value = 13.89
$ws.Range("E5").Value = 92.84
$ws.Range("F5").Value = 52.73
$ws.Range("G5").Value = 4.93
$ws.Range("H5").Value = "NW"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "08:00:19"
